# Update cryptos table rows with refreshed price/volume data (and a couple of
# re-ordered rows: Chainlink/WrappedEther swap at 15/16, Flow/EOS swap at 50/51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.870.90"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "1.752.62"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.45"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3869"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3386"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.47"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.116"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07213"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.49"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.189"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.100"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.749.60"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001060"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06606"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "79.57"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.79"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.189"
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("D23").Value = "27.867.13"
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.66"
$ws.Range("E24").Value = "  -2.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.399"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.07"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.87"
$ws.Range("E27").Value = "  -3.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.305"
$ws.Range("E28").Value = "  -4.10%  "
$ws.Range("D29").Value = "1.962.07"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.282"
$ws.Range("E30").Value = "  -10.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.05"
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.024"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.834"
$ws.Range("E33").Value = "  -4.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08795"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.16"
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.535"
$ws.Range("E36").Value = "  +2.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6537"
$ws.Range("E37").Value = "  -3.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.144"
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02275"
$ws.Range("E39").Value = "  -5.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06114"
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2108"
$ws.Range("E41").Value = "  -3.12%  "
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.012"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.71"
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.815"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6051"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.11"
$ws.Range("E48").Value = "  -3.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.995"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.113"
$ws.Range("E50").Value = "  +4.95%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.162"
$ws.Range("E51").Value = "  +1.71%  "
